$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 129
$ws.Range("H129").Value = 6017.4165
$ws.Range("I129").Value = 6199.5
$ws.Range("K129").Value = 18598.5
$ws.Range("M129").Value = -13598.5
# Row 132
$ws.Range("H132").Value = 2794.3171
$ws.Range("J132").Value = 2649.2856
$ws.Range("L132").Value = 7947.8568
$ws.Range("N132").Value = -13007.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8071.085
$ws.Range("I32").Value = 3643.5962
$ws.Range("K32").Value = 3643.5962
$ws.Range("M32").Value = -3356.5962
# Row 102
$ws.Range("H102").Value = 1649.1428
$ws.Range("I102").Value = 1507.1052
$ws.Range("J102").Value = 2998.5
$ws.Range("K102").Value = 1507.1052
$ws.Range("L102").Value = 2998.5
$ws.Range("M102").Value = 114.8948
$ws.Range("N102").Value = -6242.5
# Row 132
$ws.Range("H132").Value = 4746.483
$ws.Range("I132").Value = 2985.4167
$ws.Range("K132").Value = 8956.250100000001
$ws.Range("M132").Value = -6426.250100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 3533
$ws.Range("I107").Value = 2749.75
$ws.Range("J107").Value = 4159.6
$ws.Range("K107").Value = 2749.75
$ws.Range("L107").Value = 4159.6
$ws.Range("M107").Value = -829.75
$ws.Range("N107").Value = -7999.6
# Row 134
$ws.Range("H134").Value = 3386.375
$ws.Range("I134").Value = 2648.818
$ws.Range("K134").Value = 7946.454000000001
$ws.Range("M134").Value = -5411.454000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 519.13635
$ws.Range("I7").Value = 395.58334
$ws.Range("K7").Value = 395.58334
$ws.Range("M7").Value = -282.58334
# Row 31
$ws.Range("H31").Value = 10570.158
$ws.Range("I31").Value = 3536.5715
$ws.Range("K31").Value = 3536.5715
$ws.Range("M31").Value = -3241.5715
# Row 34
$ws.Range("H34").Value = 10570.158
$ws.Range("I34").Value = 3536.5715
$ws.Range("K34").Value = 3536.5715
$ws.Range("M34").Value = -3334.5715
# Row 107
$ws.Range("H107").Value = 6411710
$ws.Range("I107").Value = 1202.5555
$ws.Range("J107").Value = 20835352
$ws.Range("K107").Value = 1202.5555
$ws.Range("L107").Value = 20835352
$ws.Range("M107").Value = 717.4445000000001
$ws.Range("N107").Value = -20839192

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3344
$ws.Range("I3").Value = 3344
$ws.Range("K3").Value = 10032
$ws.Range("M3").Value = -9920
# Row 5
$ws.Range("H5").Value = 1797.5264
$ws.Range("I5").Value = 1753.8182
$ws.Range("J5").Value = 1857.625
$ws.Range("K5").Value = 5261.4546
$ws.Range("L5").Value = 5572.875
$ws.Range("M5").Value = -5149.4546
$ws.Range("N5").Value = -5796.875
# Row 32
$ws.Range("H32").Value = 2567.182
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 2773.9
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 8321.700000000001
$ws.Range("M32").Value = -1217
$ws.Range("N32").Value = -8887.700000000001
# Row 135
$ws.Range("H135").Value = 1797.5264
$ws.Range("I135").Value = 1753.8182
$ws.Range("J135").Value = 1857.625
$ws.Range("K135").Value = 15784.3638
$ws.Range("L135").Value = 16718.625
$ws.Range("M135").Value = -13249.3638
$ws.Range("N135").Value = -21788.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5232.6665
$ws.Range("J80").Value = 5999
$ws.Range("L80").Value = 5999
$ws.Range("N80").Value = -7995
# Row 83
$ws.Range("H83").Value = 5232.6665
$ws.Range("J83").Value = 5999
$ws.Range("L83").Value = 29995
$ws.Range("N83").Value = -39979
# Row 113
$ws.Range("H113").Value = 3313.3684
$ws.Range("I113").Value = 3155.75
$ws.Range("J113").Value = 3583.5715
$ws.Range("K113").Value = 3155.75
$ws.Range("L113").Value = 3583.5715
$ws.Range("M113").Value = -985.75
$ws.Range("N113").Value = -7923.5715
# Row 126
$ws.Range("H126").Value = 3126.4614
$ws.Range("I126").Value = 3246.8386
$ws.Range("J126").Value = 2660
$ws.Range("K126").Value = 9740.515800000001
$ws.Range("L126").Value = 7980
$ws.Range("M126").Value = -7270.515800000001
$ws.Range("N126").Value = -12920
# Row 132
$ws.Range("H132").Value = 3917.7021
$ws.Range("I132").Value = 3785.353
$ws.Range("J132").Value = 4263.846
$ws.Range("K132").Value = 11356.059
$ws.Range("L132").Value = 12791.538
$ws.Range("M132").Value = -8826.059000000001
$ws.Range("N132").Value = -17851.538

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2771.353
$ws.Range("I40").Value = 2336.6428
$ws.Range("K40").Value = 2336.6428
$ws.Range("M40").Value = -2200.6428
# Row 68
$ws.Range("H68").Value = 5017.7144
$ws.Range("I68").Value = 4849.8184
$ws.Range("J68").Value = 5633.3335
$ws.Range("K68").Value = 4849.8184
$ws.Range("L68").Value = 5633.3335
$ws.Range("M68").Value = -4100.8184
$ws.Range("N68").Value = -7131.3335
# Row 71
$ws.Range("H71").Value = 5017.7144
$ws.Range("I71").Value = 4849.8184
$ws.Range("J71").Value = 5633.3335
$ws.Range("K71").Value = 24249.092
$ws.Range("L71").Value = 28166.6675
$ws.Range("M71").Value = -20505.092
$ws.Range("N71").Value = -35654.6675
# Row 82
$ws.Range("H82").Value = 3580.4
$ws.Range("I82").Value = 1451
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 1451
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -1090
$ws.Range("N82").Value = -5722
# Row 85
$ws.Range("H85").Value = 3580.4
$ws.Range("I85").Value = 1451
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 1451
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -203
$ws.Range("N85").Value = -7496
# Row 122
$ws.Range("H122").Value = 4237.4546
$ws.Range("I122").Value = 3725.353
$ws.Range("J122").Value = 5978.6
$ws.Range("K122").Value = 11176.059
$ws.Range("L122").Value = 17935.8
$ws.Range("M122").Value = -8726.059000000001
$ws.Range("N122").Value = -22835.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 19239.125
$ws.Range("I62").Value = 14787.2
$ws.Range("K62").Value = 14787.2
$ws.Range("M62").Value = -14163.2
# Row 65
$ws.Range("H65").Value = 19239.125
$ws.Range("I65").Value = 14787.2
$ws.Range("K65").Value = 73936
$ws.Range("M65").Value = -70816
# Row 81
$ws.Range("H81").Value = 1437.5625
$ws.Range("I81").Value = 1590.1
$ws.Range("J81").Value = 1183.3334
$ws.Range("K81").Value = 3180.2
$ws.Range("L81").Value = 2366.6668
$ws.Range("M81").Value = -2119.2
$ws.Range("N81").Value = -4488.6668
# Row 84
$ws.Range("H84").Value = 1437.5625
$ws.Range("I84").Value = 1590.1
$ws.Range("J84").Value = 1183.3334
$ws.Range("K84").Value = 15901
$ws.Range("L84").Value = 11833.334
$ws.Range("M84").Value = -10597
$ws.Range("N84").Value = -22441.334
# Row 126
$ws.Range("H126").Value = 1621.9412
$ws.Range("I126").Value = 904.4
$ws.Range("K126").Value = 2713.2
$ws.Range("M126").Value = -243.1999999999998
# Row 132
$ws.Range("H132").Value = 2282.457
$ws.Range("I132").Value = 2036
$ws.Range("J132").Value = 2754.8333
$ws.Range("K132").Value = 6108
$ws.Range("L132").Value = 8264.499899999999
$ws.Range("M132").Value = -3578
$ws.Range("N132").Value = -13324.4999
